# The deck currently uses the "Integral" theme (ppt/theme/theme2.xml, the
# theme attached to the one-and-only Slide Master) which paints every slide
# with the green/yellow Integral palette. The commit swaps the presentation
# back to the plain default "Office Theme" palette (the colours that, before
# this edit, only lived - unused - in ppt/theme/theme1.xml, which the Notes
# Master points at).
#
# PowerPoint's object model edits theme colours through
# Master.Theme.ThemeColorScheme (an indexed 1-12 collection: dk1, lt1, dk2,
# lt2, accent1..accent6, hlink, folHlink) rather than by touching the raw
# XML part, so we drive the swap that way.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

# Target palette = the standard Office default theme colours.
$officeThemeColors = @(
    "000000", # 1  dk1
    "FFFFFF", # 2  lt1
    "44546A", # 3  dk2
    "E7E6E6", # 4  lt2
    "5B9BD5", # 5  accent1
    "ED7D31", # 6  accent2
    "A5A5A5", # 7  accent3
    "FFC000", # 8  accent4
    "4472C4", # 9  accent5
    "70AD47", # 10 accent6
    "0563C1", # 11 hlink
    "954F72"  # 12 folHlink
)

for ($i = 1; $i -le 12; $i++) {
    $hex = $officeThemeColors[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $colorScheme.Item($i).RGB = $r + ($g * 256) + ($b * 65536)
}
